$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 6163789632.3
$ws.Range("P2").Value = 1701792754.7
$ws.Range("Q2").Value = 950633486.34
$ws.Range("R2").Value = 91.8433678028
$ws.Range("S2").Value = 810691033.45
$ws.Range("T2").Value = 75.9977247027
$ws.Range("U2").Value = 1039919753.96
$ws.Range("V2").Value = 74.89496600210001
$ws.Range("W2").Value = 4049639671.7
$ws.Range("X2").Value = 837267976.6900001
$ws.Range("Y2").Value = 63.2906631042
$ws.Range("Z2").Value = 3700550.88
$ws.Range("AA2").Value = -87.151121527
$ws.Range("AB2").Value = 2114149960.6
$ws.Range("AC2").Value = 331.4242990546
$ws.Range("AD2").Value = 50.9404423701
$ws.Range("AE2").Value = 12.6918679527
$ws.Range("AF2").Value = 116.8579430451
$ws.Range("AG2").Value = 65.7004848199
